# The document has two BTEC/Pearson logo pictures repeated in both the
# "first page" and "default" header/footer pairs:
#   - Footers (Pearson logo, descr="...PearsonLogo.png"): rename image1.png -> image2.png
#   - Headers (BTec_Logo-Orange): rename image2.jpg -> image1.jpg
#
# These pictures are inline shapes anchored inside headers/footers, so they
# are not reachable via Document.InlineShapes - they have to be located via
# Sections(1).Headers / Sections(1).Footers and their Range.InlineShapes
# collection. InlineShape itself has no Name property in the Word object
# model, so each picture is briefly converted to a floating Shape (which
# does expose .Name), renamed, then converted back to an inline shape so
# the on-page layout/XML shape (wp:inline) is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-HeaderFooterPicture($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    [void]$shp.ConvertToInlineShape()
}

# Footers: Pearson logo -> image2.png (both the default and first-page footer)
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        Rename-HeaderFooterPicture $ftr.Range "image2.png"
    }
}

# Headers: BTEC logo -> image1.jpg (both the default and first-page header)
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        Rename-HeaderFooterPicture $hdr.Range "image1.jpg"
    }
}

Write-Output "Renamed header/footer logo pictures."
